$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 301.22223
$ws.Range("I28").Value = 243.88889
$ws.Range("J28").Value = 358.55554
$ws.Range("K28").Value = 243.88889
$ws.Range("L28").Value = 358.55554
$ws.Range("M28").Value = 241.11111
$ws.Range("N28").Value = -1328.55554
$ws.Range("H111").Value = 1189.5
$ws.Range("I111").Value = 1190.8182
$ws.Range("J111").Value = 1184.6666
$ws.Range("K111").Value = 3572.4546
$ws.Range("L111").Value = 3553.9998
$ws.Range("M111").Value = -505.4546
$ws.Range("N111").Value = -9687.9998
$ws.Range("H113").Value = 4063.4614
$ws.Range("I113").Value = 3347.2942
$ws.Range("J113").Value = 5416.222
$ws.Range("K113").Value = 3347.2942
$ws.Range("L113").Value = 5416.222
$ws.Range("M113").Value = -93.29419999999982
$ws.Range("N113").Value = -11924.222
$ws.Range("H132").Value = 4001486.2
$ws.Range("I132").Value = 4167986
$ws.Range("J132").Value = 5498
$ws.Range("K132").Value = 12503958
$ws.Range("L132").Value = 16494
$ws.Range("M132").Value = -12501428
$ws.Range("N132").Value = -21554
$ws.Range("H137").Value = 4352020.5
$ws.Range("I137").Value = 5886763
$ws.Range("K137").Value = 17660289
$ws.Range("M137").Value = -17657739
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 50002600
$ws.Range("I2").Value = 50002600
$ws.Range("K2").Value = 50002600
$ws.Range("M2").Value = -50002487
$ws.Range("H32").Value = 4957.7847
$ws.Range("I32").Value = 3675.762
$ws.Range("J32").Value = 10005.75
$ws.Range("K32").Value = 3675.762
$ws.Range("L32").Value = 10005.75
$ws.Range("M32").Value = -3388.762
$ws.Range("N32").Value = -10579.75
$ws.Range("H61").Value = 3102.2354
$ws.Range("I61").Value = 1344.8
$ws.Range("K61").Value = 1344.8
$ws.Range("M61").Value = -1132.8
$ws.Range("H74").Value = 1497.4286
$ws.Range("J74").Value = 1665.6666
$ws.Range("L74").Value = 1665.6666
$ws.Range("N74").Value = -3413.6666
$ws.Range("H77").Value = 1497.4286
$ws.Range("J77").Value = 1665.6666
$ws.Range("L77").Value = 8328.333000000001
$ws.Range("N77").Value = -17064.333
$ws.Range("H116").Value = 50002600
$ws.Range("I116").Value = 50002600
$ws.Range("K116").Value = 50002600
$ws.Range("M116").Value = -50000306
$ws.Range("H136").Value = 3102.2354
$ws.Range("I136").Value = 1344.8
$ws.Range("K136").Value = 4034.4
$ws.Range("M136").Value = -1484.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 50002600
$ws.Range("I3").Value = 50002600
$ws.Range("K3").Value = 50002600
$ws.Range("M3").Value = -50002486
$ws.Range("H17").Value = 2449.5454
$ws.Range("J17").Value = 2449.5454
$ws.Range("L17").Value = 2449.5454
$ws.Range("N17").Value = -2793.5454
$ws.Range("H134").Value = 3395
$ws.Range("I134").Value = 3090.0667
$ws.Range("J134").Value = 4538.5
$ws.Range("K134").Value = 9270.2001
$ws.Range("L134").Value = 13615.5
$ws.Range("M134").Value = -6735.2001
$ws.Range("N134").Value = -18685.5
$ws.Range("H139").Value = 29333.334
$ws.Range("J139").Value = 29333.334
$ws.Range("L139").Value = 29333.334
$ws.Range("N139").Value = -39613.334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1614731.9
$ws.Range("I31").Value = 1725403
$ws.Range("K31").Value = 1725403
$ws.Range("M31").Value = -1725108
$ws.Range("H34").Value = 1614731.9
$ws.Range("I34").Value = 1725403
$ws.Range("K34").Value = 1725403
$ws.Range("M34").Value = -1725201
$ws.Range("H58").Value = 17860178
$ws.Range("I58").Value = 1849.8572
$ws.Range("J58").Value = 35718508
$ws.Range("K58").Value = 1849.8572
$ws.Range("L58").Value = 35718508
$ws.Range("M58").Value = -1646.8572
$ws.Range("N58").Value = -35718914
$ws.Range("H99").Value = 2933
$ws.Range("J99").Value = 3279.6
$ws.Range("L99").Value = 3279.6
$ws.Range("N99").Value = -6275.6
$ws.Range("H122").Value = 2753
$ws.Range("I122").Value = 2569.9092
$ws.Range("J122").Value = 3004.75
$ws.Range("K122").Value = 7709.7276
$ws.Range("L122").Value = 9014.25
$ws.Range("M122").Value = -5259.7276
$ws.Range("N122").Value = -13914.25
$ws.Range("H126").Value = 2933
$ws.Range("J126").Value = 3279.6
$ws.Range("L126").Value = 9838.799999999999
$ws.Range("N126").Value = -14778.8
$ws.Range("H132").Value = 2952.1892
$ws.Range("I132").Value = 1696.1305
$ws.Range("J132").Value = 5015.7144
$ws.Range("K132").Value = 5088.3915
$ws.Range("L132").Value = 15047.1432
$ws.Range("M132").Value = -2558.3915
$ws.Range("N132").Value = -20107.1432
$ws.Range("H134").Value = 1790.74
$ws.Range("I134").Value = 1479.6904
$ws.Range("J134").Value = 3423.75
$ws.Range("K134").Value = 4439.0712
$ws.Range("L134").Value = 10271.25
$ws.Range("M134").Value = -1904.0712
$ws.Range("N134").Value = -15341.25
$ws.Range("H136").Value = 17860178
$ws.Range("I136").Value = 1849.8572
$ws.Range("J136").Value = 35718508
$ws.Range("K136").Value = 5549.571599999999
$ws.Range("L136").Value = 107155524
$ws.Range("M136").Value = -2999.571599999999
$ws.Range("N136").Value = -107160624
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 393059.72
$ws.Range("I126").Value = 1359.1111
$ws.Range("J126").Value = 546333.9
$ws.Range("K126").Value = 4077.3333
$ws.Range("L126").Value = 1639001.7
$ws.Range("M126").Value = -1607.3333
$ws.Range("N126").Value = -1643941.7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2483.1914
$ws.Range("I132").Value = 1583.6897
$ws.Range("J132").Value = 3932.389
$ws.Range("K132").Value = 4751.0691
$ws.Range("L132").Value = 11797.167
$ws.Range("M132").Value = -2221.0691
$ws.Range("N132").Value = -16857.167
$ws.Range("H136").Value = 2780358.5
$ws.Range("I136").Value = 4349961.5
$ws.Range("J136").Value = 3369.077
$ws.Range("K136").Value = 13049884.5
$ws.Range("L136").Value = 10107.231
$ws.Range("M136").Value = -13047334.5
$ws.Range("N136").Value = -15207.231
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 300249.44
$ws.Range("I132").Value = 402787.3
$ws.Range("K132").Value = 1208361.9
$ws.Range("M132").Value = -1205831.9
$ws.Range("H136").Value = 1651.3478
$ws.Range("I136").Value = 998.58826
$ws.Range("J136").Value = 3500.8333
$ws.Range("K136").Value = 2995.76478
$ws.Range("L136").Value = 10502.4999
$ws.Range("M136").Value = -445.76478
$ws.Range("N136").Value = -15602.4999
